# Big stimulus update:
#  - recode the `correct_ans` column (L) from single-letter codes to
#    full words: y -> left, r -> right, b -> center
#  - make the "face" image-stimulus folder/category more uniform by
#    renaming it to "book" (face//face_NN.jpg -> book//book_NN.jpg)
#    wherever it appears in the prompt/distractor file columns (A-D)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$ansMap = @{ "y" = "left"; "r" = "right"; "b" = "center" }

for ($r = 2; $r -le $rowCount; $r++) {
    # Column L = correct_ans -> expand abbreviation to full word
    $cell = $ws.Cells.Item($r, 12)
    $val = $cell.Text
    if ($ansMap.ContainsKey($val)) {
        $cell.Value = $ansMap[$val]
    }

    # Columns A-D hold the image filenames; rename the "face" stimulus
    # category to "book" (folder + filename stem) wherever it occurs.
    for ($c = 1; $c -le 4; $c++) {
        $fcell = $ws.Cells.Item($r, $c)
        $fval = $fcell.Text
        if ($fval -like "face//face_*") {
            $fcell.Value = $fval -replace "face//face_", "book//book_"
        }
    }
}
